$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)
$ws.Range("B4:I36").Copy()
$ws.Range("B3").PasteSpecial(-4163)
$ws.Rows("36:36").Delete()
